$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that was on A1 (target file no longer has <hyperlinks>)
$ws.Range("A1").Hyperlinks.Delete()

# Fix the sheet name: "Name #&2" -> "Name #2" (stray "&" removed / properly escaped)
$ws.Name = '"Name #2"'

# Replace the cell contents of the header row with simple placeholder text
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"

# Clear out the second row's text values (cells keep their style, just lose content)
$ws.Range("A2:B2").ClearContents()

# Row 2 shrinks back down to the compact row height now that the hyperlinked/long text is gone
$ws.Rows(2).RowHeight = 13.8
